$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '58.915.72'
Set-TextValue $ws.Range('E2') '  +0.30%  '
Set-TextValue $ws.Range('D3') '2.576.54'
Set-TextValue $ws.Range('E3') '  -0.66%  '
Set-TextValue $ws.Range('E4') '  -0.10%  '
Set-TextValue $ws.Range('D5') '565.39'
Set-TextValue $ws.Range('E5') '  +2.33%  '
Set-TextValue $ws.Range('D6') '143.21'
Set-TextValue $ws.Range('E6') '  -0.18%  '
Set-TextValue $ws.Range('E7') '  +0.01%  '
Set-TextValue $ws.Range('D8') '0.597'
Set-TextValue $ws.Range('E8') '  -0.02%  '
Set-TextValue $ws.Range('D9') '2.581.77'
Set-TextValue $ws.Range('E9') '  -0.68%  '
Set-TextValue $ws.Range('D10') '6.67'
Set-TextValue $ws.Range('E10') '  -1.57%  '
Set-TextValue $ws.Range('E11') '  +2.43%  '
Set-TextValue $ws.Range('D12') '0.153'
Set-TextValue $ws.Range('E12') '  +8.05%  '
Set-TextValue $ws.Range('E13') '  +1.70%  '
Set-TextValue $ws.Range('D14') '3.029.89'
Set-TextValue $ws.Range('E14') '  -0.90%  '
Set-TextValue $ws.Range('D15') '58.987.66'
Set-TextValue $ws.Range('E15') '  +0.49%  '
Set-TextValue $ws.Range('D16') '22.18'
Set-TextValue $ws.Range('E16') '  +6.66%  '
Set-TextValue $ws.Range('D17') '0.0000137'
Set-TextValue $ws.Range('E17') '  +4.02%  '
Set-TextValue $ws.Range('D18') '2.579.31'
Set-TextValue $ws.Range('E18') '  -1.50%  '
Set-TextValue $ws.Range('D19') '4.49'
Set-TextValue $ws.Range('E19') '  +1.02%  '
Set-TextValue $ws.Range('D20') '335.63'
Set-TextValue $ws.Range('E20') '  -0.30%  '
Set-TextValue $ws.Range('D21') '10.15'
Set-TextValue $ws.Range('E21') '  +1.05%  '
Set-TextValue $ws.Range('D22') '6.17'
Set-TextValue $ws.Range('E22') '  +0.57%  '
Set-TextValue $ws.Range('D23') '1.00'
Set-TextValue $ws.Range('E23') '  +0.10%  '
Set-TextValue $ws.Range('D24') '63.94'
Set-TextValue $ws.Range('E24') '  -3.81%  '
Set-TextValue $ws.Range('D25') '0.454'
Set-TextValue $ws.Range('E25') '  +6.96%  '
Set-TextValue $ws.Range('D26') '0.998'
Set-TextValue $ws.Range('E26') '  +0.22%  '
Set-TextValue $ws.Range('D27') '0.162'
Set-TextValue $ws.Range('E27') '  +2.08%  '
Set-TextValue $ws.Range('E28') '  +1.54%  '
Set-TextValue $ws.Range('D29') '0.0₃0778'
Set-TextValue $ws.Range('E29') '  +2.62%  '
Set-TextValue $ws.Range('E30') '  -0.09%  '
Set-TextValue $ws.Range('E31') '  -0.26%  '
Set-TextValue $ws.Range('D32') '6.04'
Set-TextValue $ws.Range('E32') '  +1.87%  '
Set-TextValue $ws.Range('D33') '158.24'
Set-TextValue $ws.Range('E33') '  +2.45%  '
Set-TextValue $ws.Range('D34') '18.97'
Set-TextValue $ws.Range('E34') '  +0.24%  '
Set-TextValue $ws.Range('D35') '4.01'
Set-TextValue $ws.Range('E35') '  +2.57%  '
Set-TextValue $ws.Range('B36') 'Fetch.AI'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D36') '0.879'
Set-TextValue $ws.Range('E36') '  +6.96%  '
Set-TextValue $ws.Range('D37') '0.871'
Set-TextValue $ws.Range('E37') '  -0.28%  '
Set-TextValue $ws.Range('B38') 'ImmutableX'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D38') '1.13'
Set-TextValue $ws.Range('E38') '  +1.59%  '
Set-TextValue $ws.Range('B39') 'Stacks'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D39') '1.49'
Set-TextValue $ws.Range('E39') '  +2.63%  '
Set-TextValue $ws.Range('B40') 'OKB'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D40') '36.79'
Set-TextValue $ws.Range('E40') '  -0.97%  '
Set-TextValue $ws.Range('D41') '291.49'
Set-TextValue $ws.Range('E41') '  +3.80%  '
Set-TextValue $ws.Range('D42') '3.63'
Set-TextValue $ws.Range('E42') '  +0.69%  '
Set-TextValue $ws.Range('E43') '  +0.20%  '
Set-TextValue $ws.Range('D44') '0.0969'
Set-TextValue $ws.Range('E44') '  +1.82%  '
Set-TextValue $ws.Range('D45') '0.594'
Set-TextValue $ws.Range('E45') '  -0.35%  '
Set-TextValue $ws.Range('D46') '10.62'
Set-TextValue $ws.Range('E46') '  -0.19%  '
Set-TextValue $ws.Range('E47') '  +0.33%  '
Set-TextValue $ws.Range('D48') '19.06'
Set-TextValue $ws.Range('E48') '  +1.73%  '
Set-TextValue $ws.Range('D49') '124.72'
Set-TextValue $ws.Range('E49') '  +8.67%  '
Set-TextValue $ws.Range('E50') '  +1.89%  '
Set-TextValue $ws.Range('D51') '18.51'
Set-TextValue $ws.Range('E51') '  +3.94%  '

Write-Output "Applied 98 cell updates"
